$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TESTDATA")

# New row 24: INCORRECT_SIX_DIGIT_PIN / 303030 (numeric, same style as B16/B18:B21 - left/top aligned)
$ws.Range("A24").Value = "INCORRECT_SIX_DIGIT_PIN"
$ws.Range("B24").Value = 303030
$ws.Range("B24").Style = $ws.Range("B21").Style
$ws.Range("B24").HorizontalAlignment = $ws.Range("B21").HorizontalAlignment
$ws.Range("B24").VerticalAlignment = $ws.Range("B21").VerticalAlignment

# New row 25: ALBHABETS / abcdef (text)
$ws.Range("A25").Value = "ALBHABETS"
$ws.Range("B25").Value = "abcdef"

# Extend the duplicate-value conditional formatting over the new rows
$fc = $ws.Range("A11:A23").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A11:A25"))

# Update the selection / view as captured in the saved workbook
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("C22").Select()
